# Update computed LR-pair statistics on the active sheet with the values
# produced by the refreshed TPM pipeline run ("update scripts wuth new tpm").
# Ligand/receptor expression, specificity and edge-weight columns (G, H, I, J,
# M, N, O, P, Q, R, S, T) are refreshed for data rows 2-10; all other columns
# (A-F, K, L) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08097566666666667
$ws.Range("H2").Value = 0.242927
$ws.Range("I2").Value = 0.005588990034505014
$ws.Range("J2").Value = 0.005588990034505015
$ws.Range("M2").Value = 201.098592
$ws.Range("N2").Value = 603.295776
$ws.Range("O2").Value = 0.7918622805845071
$ws.Range("P2").Value = 0.791862280584507
$ws.Range("Q2").Value = 16.284092552928
$ws.Range("R2").Value = 146.556832976352
$ws.Range("S2").Value = 0.004425710394887223
$ws.Range("T2").Value = 0.004425710394887224
$ws.Range("G3").Value = 0.08097566666666667
$ws.Range("H3").Value = 0.242927
$ws.Range("I3").Value = 0.005588990034505014
$ws.Range("J3").Value = 0.005588990034505015
$ws.Range("O3").Value = 0.1414593902976603
$ws.Range("P3").Value = 0.1414593902976603
$ws.Range("Q3").Value = 2.909013171315
$ws.Range("R3").Value = 26.181118541835
$ws.Range("S3").Value = 0.0007906151226607785
$ws.Range("T3").Value = 0.0007906151226607786
$ws.Range("G4").Value = 0.08097566666666667
$ws.Range("H4").Value = 0.242927
$ws.Range("I4").Value = 0.005588990034505014
$ws.Range("J4").Value = 0.005588990034505015
$ws.Range("O4").Value = 0.0666783291178327
$ws.Range("P4").Value = 0.06667832911783268
$ws.Range("Q4").Value = 1.371193084014444
$ws.Range("R4").Value = 12.34073775613
$ws.Range("S4").Value = 0.0003726645169570124
$ws.Range("T4").Value = 0.0003726645169570124
$ws.Range("I5").Value = 0.6976944377922635
$ws.Range("J5").Value = 0.6976944377922635
$ws.Range("M5").Value = 201.098592
$ws.Range("N5").Value = 603.295776
$ws.Range("O5").Value = 0.7918622805845071
$ws.Range("P5").Value = 0.791862280584507
$ws.Range("Q5").Value = 2032.803910640448
$ws.Range("R5").Value = 18295.23519576403
$ws.Range("S5").Value = 0.5524779086613073
$ws.Range("T5").Value = 0.5524779086613072
$ws.Range("I6").Value = 0.6976944377922635
$ws.Range("J6").Value = 0.6976944377922635
$ws.Range("O6").Value = 0.1414593902976603
$ws.Range("P6").Value = 0.1414593902976603
$ws.Range("S6").Value = 0.09869542978416246
$ws.Range("T6").Value = 0.09869542978416246
$ws.Range("I7").Value = 0.6976944377922635
$ws.Range("J7").Value = 0.6976944377922635
$ws.Range("O7").Value = 0.0666783291178327
$ws.Range("P7").Value = 0.06667832911783268
$ws.Range("S7").Value = 0.04652109934679379
$ws.Range("T7").Value = 0.04652109934679378
$ws.Range("I8").Value = 0.2967165721732315
$ws.Range("J8").Value = 0.2967165721732316
$ws.Range("M8").Value = 201.098592
$ws.Range("N8").Value = 603.295776
$ws.Range("O8").Value = 0.7918622805845071
$ws.Range("P8").Value = 0.791862280584507
$ws.Range("Q8").Value = 864.5139986699522
$ws.Range("R8").Value = 7780.625988029569
$ws.Range("S8").Value = 0.2349586615283126
$ws.Range("T8").Value = 0.2349586615283126
$ws.Range("I9").Value = 0.2967165721732315
$ws.Range("J9").Value = 0.2967165721732316
$ws.Range("O9").Value = 0.1414593902976603
$ws.Range("P9").Value = 0.1414593902976603
$ws.Range("S9").Value = 0.04197334539083704
$ws.Range("T9").Value = 0.04197334539083705
$ws.Range("I10").Value = 0.2967165721732315
$ws.Range("J10").Value = 0.2967165721732316
$ws.Range("O10").Value = 0.0666783291178327
$ws.Range("P10").Value = 0.06667832911783268
$ws.Range("S10").Value = 0.01978456525408189
$ws.Range("T10").Value = 0.01978456525408189
